$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 updates
$ws.Range("G4").Value = 2.6
$ws.Range("I4").Value = 2.8
$ws.Range("J4").Value = 3.2
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 3.5
$ws.Range("Q4").Value = 2.2
$ws.Range("V4").Value = 1.9
$ws.Range("W4").Value = 7.7
$ws.Range("Y4").Value = 9.5
$ws.Range("AA4").Value = 23
$ws.Range("AB4").Value = 32
$ws.Range("AE4").Value = 13.5
$ws.Range("AH4").Value = 14
$ws.Range("AJ4").Value = 37
$ws.Range("AK4").Value = 27
$ws.Range("AO4").Value = 14
$ws.Range("AP4").Value = 21
$ws.Range("AR4").Value = 90
$ws.Range("AW4").Value = 4.75
$ws.Range("AX4").Value = 16.5
$ws.Range("AY4").Value = 24
$ws.Range("AZ4").Value = 80
$ws.Range("BA4").Value = 120

# Row 6 updates
$ws.Range("G6").Value = 1.24
$ws.Range("H6").Value = 5.2
$ws.Range("I6").Value = 9.5
$ws.Range("J6").Value = 1.65
$ws.Range("K6").Value = 2.62
$ws.Range("L6").Value = 7.9
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 9.75
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 5
$ws.Range("Q6").Value = 1.44
$ws.Range("R6").Value = 2.6
$ws.Range("S6").Value = 1.27
$ws.Range("T6").Value = 3.45
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.87
$ws.Range("W6").Value = 9.25
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 8
$ws.Range("AB6").Value = 24
$ws.Range("AC6").Value = 9.75
$ws.Range("AD6").Value = 11.5
$ws.Range("AE6").Value = 21
$ws.Range("AF6").Value = 80
$ws.Range("AG6").Value = 32
$ws.Range("AH6").Value = 80
$ws.Range("AI6").Value = 30
$ws.Range("AJ6").Value = 300
$ws.Range("AK6").Value = 110
$ws.Range("AL6").Value = 80
$ws.Range("AM6").Value = 500
$ws.Range("AN6").Value = 3.25
$ws.Range("AO6").Value = 5.4
$ws.Range("AQ6").Value = 13
$ws.Range("AR6").Value = 35
$ws.Range("AT6").Value = 3.45
$ws.Range("AU6").Value = 8.5
$ws.Range("AV6").Value = 70
$ws.Range("AW6").Value = 10.25
$ws.Range("AX6").Value = 55
$ws.Range("AY6").Value = 45
$ws.Range("AZ6").Value = 400
$ws.Range("BA6").Value = 350

# Row 8 updates
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 2.87
$ws.Range("J8").Value = 2.77
$ws.Range("L8").Value = 3.45
$ws.Range("P8").Value = 3.8
$ws.Range("S8").Value = 1.34
$ws.Range("T8").Value = 3
$ws.Range("W8").Value = 9.75
$ws.Range("X8").Value = 13
$ws.Range("AA8").Value = 16.5
$ws.Range("AB8").Value = 22
$ws.Range("AD8").Value = 6.7
$ws.Range("AE8").Value = 11.75
$ws.Range("AF8").Value = 45
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 16.5
$ws.Range("AJ8").Value = 37
$ws.Range("AL8").Value = 27
$ws.Range("AO8").Value = 11.5
$ws.Range("AP8").Value = 17
$ws.Range("AQ8").Value = 40
$ws.Range("AR8").Value = 65
$ws.Range("AS8").Value = 175
$ws.Range("AT8").Value = 3
$ws.Range("AU8").Value = 6.5
$ws.Range("AV8").Value = 50
$ws.Range("AY8").Value = 21
$ws.Range("BA8").Value = 90
$ws.Range("BB8").Value = 250

